# Third Commit (Ruli) - katalon bio farma v.01
# Replace the ad-hoc test "Kode Vendor" codes with unified Katalon automated-test
# codes across all four data sheets, add a couple of real values that were
# missing (vendor's "Tanggal Ditetapkan"/"Tanggal Berakhir Masa Aktif" dates and
# the bank account numbers), and leave the workbook positioned on "Rekening Bank".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Vendor
# ---------------------------------------------------------------------------
$wsVendor = $wb.Worksheets.Item("Vendor")
$wsVendor.Activate() | Out-Null

$wsVendor.Range("C2").Value2 = "AutomatedTest/001"
$wsVendor.Range("C3").Value2 = "AutomatedTest/002"

# Tanggal Ditetapkan / Tanggal Berakhir Masa Aktif columns switch from free-text
# to real dates. O2/P2/O3 already carry the "yyyy-mm-dd" style (s="6"); P3 needs
# that same number format applied explicitly since it previously had none.
$wsVendor.Range("O2").Value2 = 43842
$wsVendor.Range("P2").Value2 = 45669
$wsVendor.Range("O3").Value2 = 43477
$wsVendor.Range("P3").NumberFormat = "yyyy\-mm\-dd"
$wsVendor.Range("P3").Value2 = 45700

$wsVendor.Range("Q3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Pengurus
# ---------------------------------------------------------------------------
$wsPengurus = $wb.Worksheets.Item("Pengurus")
$wsPengurus.Range("B2").Value2 = "AutomatedTest/001"
$wsPengurus.Range("B3").Value2 = "AutomatedTest/002"
$wsPengurus.Range("B2:B3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Kontak
# ---------------------------------------------------------------------------
$wsKontak = $wb.Worksheets.Item("Kontak")
$wsKontak.Range("B2").Value2 = "AutomatedTest/001"
$wsKontak.Range("B3").Value2 = "AutomatedTest/002"
$wsKontak.Range("C13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Rekening Bank
# ---------------------------------------------------------------------------
$wsRekening = $wb.Worksheets.Item("Rekening Bank")
$wsRekening.Range("B2").Value2 = "AutomatedTest/001"
$wsRekening.Range("B3").Value2 = "AutomatedTest/002"

# New bank account numbers that weren't filled in before.
$wsRekening.Range("D2").Value2 = 123456789
$wsRekening.Range("D3").Value2 = 987654321

$wsRekening.Range("D9").Select() | Out-Null

# "Rekening Bank" ends up the active/visible sheet.
$wsRekening.Activate() | Out-Null
